$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Validity End" row (A6) previously stored a raw date value (31-Dec-2025).
# Replace it with an explicit text label so it always reads clearly,
# regardless of date/number formatting.
$ws.Range("A6").Value = "Validity End: 31-Dec-2025"
